$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's record was collected and inserted at row 13 ("Hortaliza, Terminal
# La Palmera de La Serena - Cilantro" is logged daily; this insert bumps every
# following row down by one and appends the data that used to be the oldest
# tracked row as the new final row of the sheet).
$ws.Rows("13:13").Insert()

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44649
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112040
$ws.Range("G13").Value = "Cilantro"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 2300
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2400
$ws.Range("N13").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 1600
$ws.Range("Q13").Value = 1.5
$ws.Range("R13").Value = "Hortaliza"
